$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update existing orientation values (column C) for layers 1 and 2
$ws.Range("C2").Value = 90
$ws.Range("C3").Value = 90

# Add new layer rows 3, 4, 5 (rows 4, 5, 6 in the sheet)
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = 45

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = 1
$ws.Range("C5").Value = 45

$ws.Range("A6").Value = 5
$ws.Range("B6").Value = 1
$ws.Range("C6").Value = 45

# Move selection to F15, matching the final cursor position in the saved file
$ws.Range("F15").Select()
